# Update "想去人数" (want-to-go count) values in column F for the
# "展览" and "全部类型" worksheets, matching the latest scraped data.

$wb = $excel.ActiveWorkbook

$updates = @{
    "F2"  = 367
    "F3"  = 105
    "F4"  = 1569
    "F7"  = 401
    "F10" = 425
}

foreach ($sheetName in @("展览", "全部类型")) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($addr in $updates.Keys) {
        $ws.Range($addr).Value = $updates[$addr]
    }
}
